# Chaewan Woo Cover Letter (Kor) - apply edits described by the commit diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "MISS: $old"
    }
}

# 1) Title font size 36 -> 40 (half-points): sz + szCs (complex-script) for the whole title paragraph.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Font.Size = 20
$titlePara.Range.Font.SizeBi = 20

# 2) Paragraph 2: Canada study-abroad story rewrite.
Replace-Text "중학교 1학년 1학기까지 마친 후 캐나다로 유학을 떠나 Bayview Middle School을 졸업한 뒤, " `
    "중학교 1학년 1학기까지 마친 뒤, 더 넓은 시야와 다양한 경험을 쌓고자 유학을 결심하고 캐나다로 향했습니다. 낯선 환경 속에서도 빠르게 적응하며 Bayview Middle School에서 중학교 과정을 성실히 마쳤습니다. 이후 "

# 3) "9학년부터..." paragraph: add dorm-supervisor sentence.
Replace-Text "9학년부터 12학년까지 매년 성적 우수생으로 선정되었고, 그 결과 University of Toronto" `
    "9학년부터 12학년까지 매년 성적 우수생으로 선정되었고, 동시에 기숙사 내에 학생 사감으로 활동하며 공동체의 질서 유지를 돕고 후배들의 생활 적응을 지원하는 등 책임감 있는 역할도 수행했습니다. 이러한 경험을 바탕으로 University of Toronto"

# 4) "Overwatch가 세계적으로" -> "Overwatch게임이 세계적으로"
Replace-Text "University of Toronto 재학 중 Overwatch가 세계적으로 " `
    "University of Toronto 재학 중 Overwatch게임이 세계적으로 "

# 5) Club-joining sentence rewrite.
Replace-Text "고, 우연히 Overwatch 동아리를 만든 친구를 통해 게임동아리 활동에 참여하게 되었습니다." `
    "고, 직접 Overwatch게임 동아리에 참여해 전략 리더로써 활동을 주도했습니다."

# 6) "흥미를 ... 관심을 갖게 되었습니다." -> "흥미가 ... 관심으로 발전했습니다."
Replace-Text "이 시기를 계기로 게임에 대한 흥미를 단순한 취미를 넘어서서 게임의 원리와 구성, 마케팅 등에 대한 관심을 갖게 되었습니다." `
    "이 시기를 계기로 게임에 대한 흥미가 단순한 취미를 넘어서서 게임의 원리와 구성, 마케팅 등에 대한 관심으로 발전했습니다."

# 7) Military-delay / hobby paragraph big rewrite (also relocates the _GoBack bookmark).
#    A placeholder marker "@@GOBACK@@" is inserted at the exact spot the bookmark
#    must end up at; it is then located, collapsed, bookmarked, and removed.
Replace-Text "입대가 지연되었고 외부 활동도 크게 제한된 상황 속에서, 자연스럽게 취미였던 게임을 다시 접하게 되었습니다." `
    "입대가 연기되었고, 외부 활동도 크게 제한되는 답답한 시기가 이어졌습니다. 입대를 기다리며 이 시간을 의미 있게 활용하고자, 대학교 동아리에서 전략 리더로 활동 했던 경험을 바탕으로 게임 코치 또는@@GOBACK@@ 전략 분석가로 활동해보기로 결심하며, 취미였던 게임을 다시 깊이 있게 접하게 되었습니다."

# Move the _GoBack bookmark: delete the old one, then plant a new one at the marker.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$markerRange = $d.Content
$markerFound = $markerRange.Find.Execute("@@GOBACK@@")
if ($markerFound) {
    $markerRange.Text = ""
    $d.Bookmarks.Add("_GoBack", $markerRange)
} else {
    Write-Output "MISS: @@GOBACK@@ marker not found"
}

# 8) Esports achievements paragraph expansion.
Replace-Text "전 세계 랭킹 31위라는 성과를 거두었습니다. 이로 인해 여러 프로게이머 팀으로부터 입단 제의를 받았으나, 프로게이머로서의 수명이 짧다는 점, 그리고 평소에 게임을 하면서 단순히 즐기는 것에 그치지 않고 항상 새로운 게임을 접할 때마다 '어떤 시스템이나 메커니즘이 추가되면 더 좋지 않을까?'라는 고민을 해왔던 저로서는, 이러한 사고방식이 단순한 플레이어를 넘어 게임에 더 깊이 관여하고 개발에도 참여하고 싶다는 열망으로 이어졌고, 그때부터 게임 개발자의 길을 목표로 선택하게 되었습니다." `
    "전 세계 랭킹 31위라는 성과를 거두며, 일본 프로팀인 SCARZ, FENNEL의 코치와 협업하며 전략 분석가로도 활동했습니다. 이로 인해 여러 프로팀으로부터 입단 제의도 받았으나, 프로게이머로서의 수명이 짧다는 점, 그리고 단순히 게임을 즐기는 데 그치지 않고 새로운 게임을 접할 때마다 '어떤 시스템이나 메커니즘이 추가되면 더 좋지 않을까?'라는 고민을 반복해왔던 저로서는, 점차 단순한 플레이어를 넘어 게임에 더 깊이 관여하고 개발에도 참여하고 싶다는 열망을 가지게 됐습니다. 이 경험을 계기로 게임 개발자의 길을 목표로 삼게 되었습니다."

# 9) "기능 구현의 완성도..." -> add "고민하고"
Replace-Text "기능 구현의 완성도와 구조적 설계에 특히 신경 쓰며 개발 역량을 꾸준히 확장해 나갔습니다." `
    "기능 구현의 완성도와 구조적 설계에 특히  고민하고 신경 쓰며 개발 역량을 꾸준히 확장해 나갔습니다."

Write-Output "DONE"
